$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Complete import subject: update subject names in column C
$ws.Range("C2").Value = "CNTT"
$ws.Range("C3").Value = "DH"

# Fix UI admin header: update the active selection
$ws.Range("D7").Select()
